$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (state/country list) entirely.
$ws.Rows("2:4").Delete()

# Update the header row to the new column set: Alias, Name, Country ID
$ws.Range("A1").Value = "Alias"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Country ID"

# New column C needs a fitted width like the existing A/B columns.
$ws.Columns("C:C").AutoFit()

